$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Selection moves from K9 to H15.
$ws.Range("H15").Select()

# M12 switches from the shared string "бендер" to a new "S бендер" string
# ("Revert 'Revert \"s bendes\"'").
$ws.Range("M12").Value = "S бендер"

# Row 14 gains a new B14 cell (same "4.5-20" value used in the analogous
# B12 cell for this instrument block).
$ws.Range("B14").Value = "4.5-20"

# The whole row 14 (A:M) is highlighted with the workbook's existing
# yellow fill, matching the formatting already used for rows 6/8/10 etc.
# This also creates the (previously empty/absent) I14 and J14 cells.
$ws.Range("A14:M14").Interior.Color = 65535
